$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.494.64'
$ws.Range("E2").Value = '  -0.16%  '

$ws.Range("D3").Value = '2.326.53'
$ws.Range("E3").Value = '  +0.12%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '512.39'
$ws.Range("E5").Value = '  -1.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.83'
$ws.Range("E6").Value = '  -2.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.533'
$ws.Range("E8").Value = '  -0.62%  '

$ws.Range("E9").Value = '  -3.19%  '

$ws.Range("E10").Value = '  -0.21%  '

$ws.Range("E11").Value = '  -0.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.337'
$ws.Range("E12").Value = '  -1.27%  '

$ws.Range("D13").Value = '2.741.64'
$ws.Range("E13").Value = '  +0.04%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.51'
$ws.Range("E14").Value = '  -0.97%  '

$ws.Range("D15").Value = '56.470.87'
$ws.Range("E15").Value = '  -0.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000133'
$ws.Range("E16").Value = '  -1.47%  '

$ws.Range("D17").Value = '2.328.44'
$ws.Range("E17").Value = '  -0.71%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.44'
$ws.Range("E18").Value = '  -0.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '326.91'
$ws.Range("E19").Value = '  +1.16%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.13'
$ws.Range("E20").Value = '  -2.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.76'
$ws.Range("E21").Value = '  +3.61%  '

$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.66'
$ws.Range("E23").Value = '  +1.83%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.65'
$ws.Range("E24").Value = '  +9.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.163'
$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.31'
$ws.Range("E27").Value = '  +3.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '167.35'
$ws.Range("E28").Value = '  -1.54%  '

$ws.Range("E29").Value = '  -2.83%  '

$ws.Range("D30").Value = '0.0₃0718'
$ws.Range("E30").Value = '  -3.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.12'
$ws.Range("E31").Value = '  -0.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.37'
$ws.Range("E32").Value = '  +0.74%  '

$ws.Range("E33").Value = '  +0.00%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").Value = '  +0.52%  '

$ws.Range("E35").Value = '  +1.33%  '

$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.91'
$ws.Range("E36").Value = '  -2.29%  '

$ws.Range("B37").Value = 'SuiNetwork'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.887'
$ws.Range("E37").Value = '  -3.77%  '

$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.55'
$ws.Range("E38").Value = '  +0.19%  '

$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '38.41'
$ws.Range("E39").Value = '  +1.34%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '150.30'
$ws.Range("E40").Value = '  +9.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.373'
$ws.Range("E41").Value = '  -1.65%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.56'
$ws.Range("E42").Value = '  -0.41%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '279.27'
$ws.Range("E43").Value = '  +0.88%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.04'
$ws.Range("E44").Value = '  -1.76%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0929'
$ws.Range("E45").Value = '  -0.33%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0494'
$ws.Range("E46").Value = '  -1.78%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.558'
$ws.Range("E47").Value = '  -0.72%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.12'
$ws.Range("E48").Value = '  +2.26%  '

$ws.Range("B49").Value = 'Polygon'
$ws.Range("C49").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.379'
$ws.Range("E49").Value = '  -0.05%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0214'
$ws.Range("E50").Value = '  -1.41%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.07'
$ws.Range("E51").Value = '  +1.77%  '
